$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7: fix capitalisation of the comment in H7 ---
$ws.Range("H7").Value = "Start Iteration 4"

# --- New time-log entries for the second working day (Iteration 4 continued) ---
# Rows 10 and 11 currently have no date formatting on column C, so copy the
# existing m/d/yyyy date format (as used by C7/C8/C9) down into C10:C11 first.
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C10:C11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 10
$ws.Range("B10").Value = "Coding"
$ws.Range("C10").Value = "4/28/2019"
$ws.Range("D10").Value = "2:20pm"
$ws.Range("F10").Value = "3:42pm"
$ws.Range("H10").Value = "Start coding appendix 2.4"

# Row 11
$ws.Range("B11").Value = "Testing"
$ws.Range("C11").Value = "4/28/2019"
$ws.Range("D11").Value = "8:12pm"
$ws.Range("F11").Value = "8:47pm"
$ws.Range("H11").Value = "Write the tests"

# Row 12 (C12 already carries the date format from the template)
$ws.Range("B12").Value = "Coding"
$ws.Range("C12").Value = "4/28/2019"
$ws.Range("D12").Value = "8:48pm"
$ws.Range("F12").Value = "9:00pm"

# Row 13 (C13 already carries the date format from the template)
$ws.Range("B13").Value = "Testing"
$ws.Range("C13").Value = "4/28/2019"
$ws.Range("D13").Value = "8:32pm"
$ws.Range("F13").Value = "9:53pm"
$ws.Range("H13").Value = "More testing and writing test automation"

# --- Update the active selection to reflect where the author ended up ---
$ws.Range("H13").Select()
